$d = $word.ActiveDocument

$replacements = @(
    @{old='150×7=1050'; new='981×8=7848'},
    @{old='715×6=4290'; new='353×4=1412'},
    @{old='518×7=3626'; new='733×8=5864'},
    @{old='289×2=578';  new='241×4=964'},
    @{old='456×9=4104'; new='280×3=840'},
    @{old='329×7=2303'; new='626×8=5008'},
    @{old='326×9=2934'; new='408×5=2040'},
    @{old='481×7=3367'; new='224×3=672'},
    @{old='398×9=3582'; new='674×7=4718'},
    @{old='725×4=2900'; new='906×5=4530'},
    @{old='168×8=1344'; new='658×2=1316'},
    @{old='738×2=1476'; new='193×4=772'},
    @{old='827×9=7443'; new='733×4=2932'},
    @{old='698×6=4188'; new='512×2=1024'},
    @{old='287×6=1722'; new='298×9=2682'},
    @{old='957×8=7656'; new='802×8=6416'},
    @{old='265×9=2385'; new='678×8=5424'},
    @{old='898×5=4490'; new='602×6=3612'},
    @{old='606×8=4848'; new='457×7=3199'},
    @{old='911×6=5466'; new='221×7=1547'},
    @{old='466×6=2796'; new='328×8=2624'},
    @{old='321×8=2568'; new='102×2=204'},
    @{old='498×4=1992'; new='642×4=2568'},
    @{old='842×7=5894'; new='705×9=6345'},
    @{old='856×8=6848'; new='815×5=4075'}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2) | Out-Null
}
